$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells whose new values could be
# misinterpreted as numbers (e.g. "1.00", "3.70"), so they stay literal text
# exactly like the source inline strings.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D25", "D27", "D28", "D29", "D31", "D33", "D38", "D39", "D42", "D43", "D46", "D48", "D49", "D50")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated cell values.
$ws.Range('D2').Value = '67.957.88'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '2.634.64'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '597.56'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '153.81'
$ws.Range('E6').Value = '  +1.02%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').Value = '2.633.44'
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('D10').Value = '0.135'
$ws.Range('E10').Value = '  +10.37%  '
$ws.Range('E11').Value = '  -0.57%  '
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('E15').Value = '  +3.94%  '
$ws.Range('D16').Value = '3.114.41'
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('D17').Value = '67.762.35'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').Value = '2.635.48'
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('D19').Value = '375.51'
$ws.Range('E19').Value = '  +3.49%  '
$ws.Range('D20').Value = '11.39'
$ws.Range('E20').Value = '  +2.27%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').Value = '4.26'
$ws.Range('D23').Value = '4.81'
$ws.Range('E23').Value = '  -1.32%  '
$ws.Range('E24').Value = '  -2.59%  '
$ws.Range('D25').Value = '72.35'
$ws.Range('E25').Value = '  +2.24%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '9.96'
$ws.Range('E27').Value = '  -0.47%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.764.00'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0000104'
$ws.Range('E29').Value = '  +2.45%  '
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('D31').Value = '576.49'
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').Value = '7.87'
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  -1.36%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = '158.34'
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').Value = '19.17'
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('E40').Value = '  +5.67%  '
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').Value = '5.34'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').Value = '2.62'
$ws.Range('E43').Value = '  +3.00%  '
$ws.Range('E44').Value = '  +13.11%  '
$ws.Range('E45').Value = '  +4.81%  '
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').Value = '155.79'
$ws.Range('D49').Value = '3.70'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').Value = '22.01'
$ws.Range('E50').Value = '  +7.02%  '
$ws.Range('E51').Value = '  -1.38%  '
